# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = -7
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -15
$ws.Range("F15").Value = -1
